# Generate Report for handback
# The localization file "e900ff5d-9406-4620-b0b0-4ecc073d7efd.md" has now been
# handed back (in sync with en-US) for both the zh-cn and de-de locales.
# Update the Overview sheet and each locale sheet accordingly.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $statusHandedBack
$wsZhCn.Range("G3").Value = "2016-01-17 07:36:09"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $statusHandedBack
$wsDeDe.Range("G3").Value = "2016-01-17 07:36:26"
